# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same per-fund-holding layout) to create
#    the new "2022-Q1" sheet, positioned right after "2021-Q4" and before
#    "总计". Trim it down to a single data row and overwrite that row with
#    the 2022-Q1 fund-holding figures.
# 2. Insert a new top data row in "总计" for the "2022-Q1" summary line,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Force the cell to store $val as text (matches the source workbook's
    # t="inlineStr" cells) instead of letting Excel auto-convert numeric-
    # looking strings ("161224", "0.77", ...) into real numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet from "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# "2021-Q4" has two holding rows; 2022-Q1 only needs one, so drop row 3.
$newSheet.Rows.Item(3).Delete()

Set-TextValue $newSheet.Range("B2") "161224"
Set-TextValue $newSheet.Range("C2") "国投瑞银新丝路灵活配置混合(LOF)"
Set-TextValue $newSheet.Range("D2") "0.77"
Set-TextValue $newSheet.Range("E2") "94.48"
Set-TextValue $newSheet.Range("F2") "3.95"
Set-TextValue $newSheet.Range("G2") "0.0304"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q1" row at the top of "总计"
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()

# Restore the row-label formatting (style index carried by column A) by
# copying it down from the row that just got pushed to A3.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2:D2").Style = "Normal"
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.03

# The A column is a plain 0-based row counter (not a formula), so bump the
# rows that were pushed down by the insert.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3

# Restore the originally-active tab ("2021-Q2") so this edit doesn't shift
# the workbook's selected sheet as a side effect.
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "2022-Q1 sheet added and 总计 updated"
